# Add new diary rows for 2.7, 2.8 & 2.9 (done earlier) to the "part2" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("part2")

# Row 5: ex 2.7
$ws.Cells.Item(5, 1).Value = 211109
$ws.Cells.Item(5, 2).Value = 30
$ws.Cells.Item(5, 3).Value = "ex 2.7"

# Row 6: ex 2.8 & 2.9
$ws.Cells.Item(6, 1).Value = 211117
$ws.Cells.Item(6, 2).Formula = "=30+10"
$ws.Cells.Item(6, 3).Value = "ex 2.8 & 2.9"

# Update the selection to match the author's final cursor position.
$ws.Activate()
$ws.Range("C11").Select()
